# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de) this "handback" run:
#   - fills in the "Latest Target File" / "Latest Handback File" columns
#     (F/G) with hyperlinks, mirroring the existing handoff hyperlinks
#   - stamps "Latest Handback DateTime" (H) with the handback timestamp
#   - flips the Status column from "Ready for handoff" to
#     "Handed back: in sync with en-US"
# The Status text lives in a shared string that both the Overview rollup
# sheet and the two locale sheets all point at, so every cell showing it
# gets updated to stay in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$sourceRepoBase = "https://github.com/OpenLocalizationTest/oltest/blob/c417183c72e021f2fc3b50ea72a299ac5d51b747/e2e"

# --- Overview sheet: Status rollup columns (B = zh-cn, C = de-de) -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- per-locale handback details ---------------------------------------
$locales = @(
    @{
        Sheet        = "zh-cn"
        HandbackTime = "2016-03-18 14:27:14"
        XlfUrl       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15aebd9b499008c264a0395d8667bdc1544045f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        XlfName      = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
    },
    @{
        Sheet        = "de-de"
        HandbackTime = "2016-03-18 14:27:22"
        XlfUrl       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ae16b14e0b30230693bde5f0e6559e5e5b2829d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        XlfName      = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    for ($row = 2; $row -le 3; $row++) {

        # Status (column C)
        $ws.Range("C$row").Value = $newStatus

        # Latest Target File (F) - the en-US source, now confirmed in sync
        $ws.Hyperlinks.Add($ws.Range("F$row"), "$sourceRepoBase/a.md", "", "", "a.md") | Out-Null

        # Latest Handback File (G) - the translated xlf handed back
        $ws.Hyperlinks.Add($ws.Range("G$row"), $locale.XlfUrl, "", "", $locale.XlfName) | Out-Null

        # Latest Handback DateTime (H)
        $ws.Range("H$row").Value = $locale.HandbackTime
    }
}
